# Preparation for transport:
#  - Independent num/denum conversion
#  - Added some passenger convs
#  - CAP2ACT is now entity dependent
#
# The underlying edit: a new "capacity_to_activity" parameter row is
# inserted right after the "buildrate" row (row 9) and before the old
# "co2_factor" row, i.e. at row 10 of Sheet1. Everything below shifts
# down by one row (dimension, autofilter, filter-database defined name
# and the selection all need to follow).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new blank row at row 10 - pushes every row from 10 on
#    down by one (co2_factor, learning_rate, ... down to the trailing
#    anchor row all shift from row N to row N+1).
$ws.Rows.Item(10).Insert()

# 2. Populate the freshly inserted row 10 with the new
#    "capacity_to_activity" parameter entry for conv_chp_nuclear.
$ws.Range("A10").Value = "CHE"
$ws.Range("B10").Value = "conv_chp_nuclear"
$ws.Range("C10").Value = "capacity_to_activity"
$ws.Range("D10").Value = "constant"
$ws.Range("G10").Value = 0.001
$ws.Range("H10").Value = "GW/TWh"

# 3. The used range grew by one row (was A1:L429 -> now A1:L430); the
#    autofilter and the hidden _FilterDatabase defined name both
#    covered A5:L849 and must now cover A5:L850.
$ws.AutoFilterMode = $false
$ws.Range("A5:L850").AutoFilter()

$filterDbName = $wb.Names.Item(1)
$filterDbName.RefersTo = "=Sheet1!`$A`$5:`$L`$850"

# 4. Restore the cursor/selection to what the author left it at:
#    B9:B10 with B9 active.
$ws.Range("B9:B10").Select()
